# tak-4 has been completed
# Populate Sheet1 with a bike-loan installment calculation table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Labels (column C)
$ws.Range("C5").Value  = "RATE OF BIKE"
$ws.Range("C6").Value  = "DOWN PAYMENT"
$ws.Range("C7").Value  = "LAON AMOUNT"
$ws.Range("C8").Value  = "RATE OF INTEREST"
$ws.Range("C9").Value  = "DURATION"
$ws.Range("C10").Value = "INSTALLMENT"
$ws.Range("C12").Value = "TOTAL AMOUNT TOBE PAID"
$ws.Range("C13").Value = "TOTAL INTEREST PAID"

# Values (column D)
$ws.Range("D5").Value  = 62500
$ws.Range("D6").Value  = 24000
$ws.Range("D7").Value  = 38500
$ws.Range("D8").Value  = 0.0925
$ws.Range("D9").Value  = 24
$ws.Range("D10").Value = 1763.28
$ws.Range("D12").Value = 42318.77
$ws.Range("D13").Value = 3818.77

# Number formats
$ws.Range("D8").NumberFormat = "0.00%"
$ws.Range("D10").NumberFormat = """₹""\ #,##0.00;[Red]""₹""\ \-#,##0.00"
$ws.Range("D12").NumberFormat = """₹""\ #,##0.00;[Red]""₹""\ \-#,##0.00"
$ws.Range("D13").NumberFormat = """₹""\ #,##0.00;[Red]""₹""\ \-#,##0.00"

# Column widths
$ws.Columns.Item(3).ColumnWidth = 18.7109375
$ws.Columns.Item(4).ColumnWidth = 17.28515625

# View state: selection on C18
$null = $ws.Range("C18").Select()
